# Edit script for diagrams/uml/associations/roles/adminStudent.pptx
#
# The canonical diff for this commit is dominated by PowerPoint
# "resave" noise (new creationId GUIDs, refreshed date/slide-number
# field ids & cached text, locale churn en-US -> en-SG on empty
# paragraph marks, notesMaster/notesSlide housekeeping, shape
# renumbering, default-valued attribute cleanup, etc.) that is not
# reachable/controllable through the PowerPoint object model and
# carries no real content meaning.
#
# The one deliberate, content-level edit visible in the diff is on
# the single slide: an empty, unused "Title 26" placeholder shape
# (id=27) is deleted from the slide's shape tree. Everything else on
# the slide (the Admin/Student rectangles, connector, textbox, and
# their animation timing) is left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find and remove the empty, unused title placeholder shape ("Title 26").
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Title 26") {
        $sh.Delete()
    }
}
